$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the numeric-looking Price/Volume columns so that
# values such as "256.30" or "4.30%" are preserved exactly as strings
# (matching the original inlineStr / shared-string cell contents) instead
# of being reinterpreted by Excel as numbers/percentages.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = "256.30"
$ws.Range("E2").Value = "4.30%"
$ws.Range("D3").Value = "28.18"
$ws.Range("E3").Value = "-3.91%"
$ws.Range("D4").Value = "5.293"
$ws.Range("E4").Value = "2.50%"
$ws.Range("D5").Value = "0.05805"
$ws.Range("E5").Value = "0.61%"
$ws.Range("D6").Value = "6.696"
$ws.Range("E6").Value = "1.37%"
$ws.Range("D7").Value = "3.224"
$ws.Range("E7").Value = "1.95%"
$ws.Range("D8").Value = "0.8718"
$ws.Range("E8").Value = "1.52%"
$ws.Range("D9").Value = "0.9034"
$ws.Range("E9").Value = "5.13%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1411"
$ws.Range("E10").Value = "3.44%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07180"
$ws.Range("E11").Value = "2.40%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03159"
$ws.Range("E12").Value = "5.40%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09237"
$ws.Range("E13").Value = "-1.35%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001538"
$ws.Range("E14").Value = "0.02%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0006073"
$ws.Range("E15").Value = "1.42%"
$ws.Range("D16").Value = "0.005889"
$ws.Range("E16").Value = "-2.86%"
$ws.Range("D17").Value = "3.507"
$ws.Range("E17").Value = "0.25%"
$ws.Range("E18").Value = "4.67%"
$ws.Range("E19").Value = "-2.32%"
$ws.Range("D20").Value = "0.03424"
$ws.Range("E20").Value = "3.94%"
$ws.Range("D21").Value = "0.1312"
$ws.Range("E21").Value = "2.38%"
$ws.Range("D22").Value = "3.523"
$ws.Range("E22").Value = "10.78%"
$ws.Range("D23").Value = "0.04154"
$ws.Range("E23").Value = "0.11%"
$ws.Range("E25").Value = "-0.25%"
$ws.Range("D26").Value = "0.004974"
$ws.Range("E26").Value = "20.32%"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").Value = "-0.91%"
$ws.Range("D28").Value = "0.0001936"
$ws.Range("E28").Value = "33.64%"
$ws.Range("D40").Value = "0.03872"
$ws.Range("E40").Value = "3.79%"
$ws.Range("D41").Value = "0.005759"
$ws.Range("E41").Value = "64.39%"
$ws.Range("E42").Value = "2.48%"
$ws.Range("D43").Value = "0.002313"
$ws.Range("E43").Value = "-5.18%"
$ws.Range("D44").Value = "0.01076"
$ws.Range("E44").Value = "27.46%"
$ws.Range("D45").Value = "0.00005264"
$ws.Range("E45").Value = "-0.41%"
$ws.Range("E46").Value = "-0.09%"
$ws.Range("D47").Value = "0.08490"
$ws.Range("E47").Value = "46.42%"
$ws.Range("D48").Value = "0.002181"
$ws.Range("E48").Value = "-4.02%"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.09%"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.09%"

# Restore the default (unstyled) cell style so the cells keep matching
# the original workbook formatting (no explicit style index).
$ws.Range("D2:E50").Style = "Normal"
